$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer samples (rows 961-1000, columns A:F) appended below the
# existing data that ended at row 960 ("V2" upload).
$data = @(
@(0.0040000000000000001,0.46200000000000002,0.81399999999999995,267.08999999999997,22.949000000000002,-46.57),
@(0.057000000000000002,0.23200000000000001,0.77100000000000002,132.32400000000001,2.6859999999999999,-29.297000000000001),
@(0.033000000000000002,0.19500000000000001,0.68100000000000005,-64.209000000000003,15.259,-8.6669999999999998),
@(0.0080000000000000002,0.45400000000000001,0.82299999999999995,-329.89499999999998,2.9910000000000001,44.25),
@(0.11600000000000001,0.626,1.1499999999999999,-249.268,-3.2349999999999999,73.058999999999997),
@(-0.025000000000000001,0.64300000000000002,0.873,166.87,-3.54,-49.866),
@(-0.070999999999999994,0.495,0.81499999999999995,265.25900000000001,4.944,-47.667999999999999),
@(-0.012,0.28799999999999998,0.72299999999999998,130.554,-27.161000000000001,-16.356999999999999),
@(-0.13900000000000001,0.34999999999999998,0.73999999999999999,-179.01599999999999,8.6669999999999998,13.855),
@(-0.14699999999999999,0.49399999999999999,1.0049999999999999,-285.88900000000001,-18.004999999999999,68.787000000000006),
@(-0.20799999999999999,0.53300000000000003,1.101,211.91399999999999,-18.494,-14.587),
@(-0.22700000000000001,0.438,0.82199999999999995,257.38499999999999,24.719000000000001,-22.399999999999999),
@(-0.153,0.17000000000000001,0.83599999999999997,202.02600000000001,-19.835999999999999,-6.7750000000000004),
@(-0.078,0.14199999999999999,0.72399999999999998,-48.279000000000003,18.187999999999999,14.099),
@(-0.193,0.36799999999999999,0.75800000000000001,-337.58499999999998,9.3989999999999991,23.376000000000001),
@(-0.055,0.45000000000000001,1.0589999999999999,-206.17699999999999,-20.751999999999999,33.752000000000002),
@(-0.021999999999999999,0.59399999999999997,0.878,44.25,25.818000000000001,-19.774999999999999),
@(0,0.496,0.84299999999999997,163.81800000000001,5.0049999999999999,-33.020000000000003),
@(0.050000000000000003,0.33900000000000002,0.73599999999999999,195.435,-2.4409999999999998,-29.358000000000001),
@(0.032000000000000001,0.29499999999999998,0.77600000000000002,-39.063000000000002,-14.709,-1.343),
@(0.050000000000000003,0.46700000000000003,0.93799999999999994,-277.52699999999999,-33.752000000000002,25.696000000000002),
@(0.075999999999999998,0.51600000000000001,1.3460000000000001,-63.904000000000003,-51.758000000000003,1.038),
@(0.059999999999999998,0.47099999999999997,0.89100000000000001,328.18599999999998,24.536000000000001,-38.573999999999998),
@(0.090999999999999998,0.191,0.78500000000000003,258.60599999999999,20.568999999999999,-30.762),
@(0.125,0.085000000000000006,0.66000000000000003,-27.893000000000001,-21.606000000000002,-13.305999999999999),
@(0.012999999999999999,0.32100000000000001,1.002,-334.71699999999998,-18.187999999999999,-0.54900000000000004),
@(0.017999999999999999,0.56200000000000006,0.93999999999999995,-291.68700000000001,-27.832000000000001,15.503),
@(0.079000000000000001,0.56499999999999995,1.097,15.076000000000001,4.8220000000000001,-12.817),
@(0.088999999999999996,0.56299999999999994,0.81399999999999995,345.642,30.762,-10.497999999999999),
@(0.18099999999999999,0.245,0.76000000000000001,257.935,2.5019999999999998,1.831),
@(0.20599999999999999,0.154,0.63800000000000001,-31.981999999999999,-3.2959999999999998,15.015000000000001),
@(0.017999999999999999,0.34999999999999998,0.98299999999999998,-319.03100000000001,-17.821999999999999,14.954000000000001),
@(0.037999999999999999,0.501,1.0760000000000001,-231.506,-32.103999999999999,28.076000000000001),
@(-0.083000000000000004,0.47199999999999998,1.089,233.459,24.902000000000001,-15.259),
@(-0.024,0.42299999999999999,0.81699999999999995,357.42200000000003,40.161000000000001,35.094999999999999),
@(0.074999999999999997,0.23000000000000001,0.81000000000000005,190.73500000000001,-7.9349999999999996,17.760999999999999),
@(0.025000000000000001,0.097000000000000003,0.68799999999999994,-123.474,16.356999999999999,-28.381),
@(-0.109,0.28999999999999998,1.0569999999999999,-377.93000000000001,-16.785,-3.9060000000000001),
@(-0.040000000000000001,0.69999999999999996,0.95599999999999996,-350.76900000000001,41.381999999999998,42.296999999999997),
@(-0.050999999999999997,0.65100000000000002,1.006,131.65299999999999,11.901999999999999,-27.344000000000001)
)

$startRow = 961
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
    }
}

# Match the saved view state from the diff: scrolled so row 977 is at the
# top, with H989 as the active/selected cell.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 977
$win.ScrollColumn = 1
$ws.Range("H989").Select()
